$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 14-25 hold the "Nordeste" region series. The diff shifts each row's
# date label forward by one year and replaces the value with the figure
# from the refreshed source data.
#
# NumberFormat is forced to Text ("@") before writing the date-like string
# so Excel doesn't auto-convert "01/12/2010" into a date serial; the style
# is then reset back to "Normal" so the cell keeps the workbook's original
# (unstyled) look, matching the source formatting.

$updates = @(
    @{ Row = 14; Date = "01/12/2010"; Value = 10.61906913691821 },
    @{ Row = 15; Date = "01/12/2011"; Value = 5.47385620915033 },
    @{ Row = 16; Date = "01/12/2012"; Value = 8.068680609346778 },
    @{ Row = 17; Date = "01/12/2013"; Value = 7.322900489786166 },
    @{ Row = 18; Date = "01/12/2014"; Value = -0.7902938557435513 },
    @{ Row = 19; Date = "01/12/2015"; Value = -9.39077751598788 },
    @{ Row = 20; Date = "01/12/2016"; Value = -5.819712729073789 },
    @{ Row = 21; Date = "01/12/2017"; Value = 3.15540362871416 },
    @{ Row = 22; Date = "01/12/2018"; Value = 0.2549069589599773 },
    @{ Row = 23; Date = "01/12/2019"; Value = 3.114670734808023 },
    @{ Row = 24; Date = "01/12/2020"; Value = 1.356182961410446 },
    @{ Row = 25; Date = "01/12/2021"; Value = -8.307991728500197 }
)

foreach ($u in $updates) {
    $cCell = $ws.Cells.Item($u.Row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $u.Date
    $cCell.Style = "Normal"

    $ws.Cells.Item($u.Row, 4).Value = $u.Value
}
